$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.027.05"
$ws.Range("E2").Value = "  +0.60%  "

$ws.Range("D3").Value = "1.645.33"
$ws.Range("E3").Value = "  +0.98%  "

$ws.Range("E4").Value = "  +0.66%  "

$ws.Range("D5").Value = "'216.43"
$ws.Range("E5").Value = "  +0.85%  "

$ws.Range("D6").Value = "'0.507"
$ws.Range("E6").Value = "  +1.11%  "

$ws.Range("E7").Value = "  +0.62%  "

$ws.Range("D8").Value = "'0.256"
$ws.Range("E8").Value = "  +0.64%  "

$ws.Range("D9").Value = "'0.0640"
$ws.Range("E9").Value = "  +1.41%  "

$ws.Range("D10").Value = "'19.66"
$ws.Range("E10").Value = "  +0.03%  "

$ws.Range("D12").Value = "1.874.72"

$ws.Range("E13").Value = "  +1.22%  "

$ws.Range("D14").Value = "1.649.51"
$ws.Range("E14").Value = "  +2.41%  "

$ws.Range("E15").Value = "  +0.27%  "

$ws.Range("D16").Value = "0.0₃0765"
$ws.Range("E16").Value = "  +0.79%  "

$ws.Range("D17").Value = "'63.22"
$ws.Range("E17").Value = "  +0.77%  "

$ws.Range("D18").Value = "26.046.59"
$ws.Range("E18").Value = "  +0.69%  "

$ws.Range("E19").Value = "  +0.68%  "

$ws.Range("D20").Value = "'193.09"
$ws.Range("E20").Value = "  +0.21%  "

$ws.Range("E21").Value = "  -0.61%  "

$ws.Range("D22").Value = "'9.95"
$ws.Range("E22").Value = "  -0.09%  "

$ws.Range("E23").Value = "  +0.08%  "

$ws.Range("B24").Value = "Stellar"
$ws.Range("C24").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D24").Value = "'0.132"
$ws.Range("E24").Value = "  +5.30%  "

$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "'1.80"
$ws.Range("E25").Value = "  +1.13%  "

$ws.Range("D26").Value = "'144.31"
$ws.Range("E26").Value = "  +0.78%  "

$ws.Range("E27").Value = "  +0.77%  "

$ws.Range("D28").Value = "'6.92"
$ws.Range("E28").Value = "  +0.83%  "

$ws.Range("E29").Value = "  +0.72%  "

$ws.Range("E30").Value = "  +1.15%  "

$ws.Range("E31").Value = "  +0.43%  "

$ws.Range("E32").Value = "  -0.59%  "

$ws.Range("D33").Value = "'3.25"
$ws.Range("E33").Value = "  +1.22%  "

$ws.Range("E34").Value = "  -3.09%  "

$ws.Range("E35").Value = "  +2.18%  "

$ws.Range("D36").Value = "'0.905"
$ws.Range("E36").Value = "  +0.44%  "

$ws.Range("D37").Value = "1.129.16"
$ws.Range("E37").Value = "  -0.66%  "

$ws.Range("E38").Value = "  -1.46%  "

$ws.Range("E39").Value = "  +0.68%  "

$ws.Range("E40").Value = "  +0.71%  "

$ws.Range("D41").Value = "'5.51"
$ws.Range("E41").Value = "  +0.87%  "

$ws.Range("D42").Value = "'99.37"
$ws.Range("E42").Value = "  +0.30%  "

$ws.Range("D43").Value = "'0.797"
$ws.Range("E43").Value = "  -0.30%  "

$ws.Range("D44").Value = "1.783.99"
$ws.Range("E44").Value = "  +1.09%  "

$ws.Range("E45").Value = "  +4.02%  "

$ws.Range("D46").Value = "'56.62"
$ws.Range("E46").Value = "  +0.85%  "

$ws.Range("D47").Value = "'0.0530"
$ws.Range("E47").Value = "  +0.17%  "

$ws.Range("D48").Value = "'1.45"
$ws.Range("E48").Value = "  +0.08%  "

$ws.Range("D49").Value = "'7.70"
$ws.Range("E49").Value = "  +1.07%  "

$ws.Range("E50").Value = "  +0.27%  "

$ws.Range("D51").Value = "'0.0957"
$ws.Range("E51").Value = "  -0.52%  "
